# Append a new day's snapshot pair (columns AV/AW) to the rolling attack-log
# table, mirroring the existing AT/AU pair, and promote the now-stale AU
# column (previously inline-string numbers) to real numeric cells.
#
# Layout recap (per row):
#   AT = styled numeric "activity" marker (fill colour encodes status)
#   AU = plain numeric "total" value, historically stored as inline text
# New columns:
#   AV = copy of AT (same style + value)    -> new "activity" marker
#   AW = copy of AU's old (text) value      -> new "total" snapshot
#   AU = rewritten in place as a real number (was inline text)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Discover the current extent: last data row, and the existing "AU" (last)
# column, so AT/AU/AV/AW resolve correctly even if the sheet shape shifts.
$used = $ws.UsedRange
$lastRow = $used.Rows.Count()
$colAU = $used.Columns.Count()
$colAT = $colAU - 1
$colAV = $colAU + 1
$colAW = $colAU + 2

# Row 1: header labels for the new date pair.
$atHeader = $ws.Cells.Item(1, $colAT)
$avHeader = $ws.Cells.Item(1, $colAV)
$awHeader = $ws.Cells.Item(1, $colAW)
$atHeader.Copy($avHeader)
$atHeader.Copy($awHeader)
$avHeader.Value = "06-11_A"
$awHeader.Value = "06-11_0"

for ($r = 2; $r -le $lastRow; $r++) {
    $atCell = $ws.Cells.Item($r, $colAT)
    $auCell = $ws.Cells.Item($r, $colAU)
    $avCell = $ws.Cells.Item($r, $colAV)
    $awCell = $ws.Cells.Item($r, $colAW)

    # New AV mirrors AT (style + value, numeric or blank).
    $atCell.Copy($avCell)

    # New AW takes AU's current (pre-conversion) value/type verbatim.
    $auCell.Copy($awCell)

    # Promote AU from inline-string to a genuine number, in place, only
    # when it actually holds text (skip truly blank cells).
    $auText = $auCell.Text()
    if ($auText -ne "") {
        $auCell.Value = $auCell.Value()
    }
}
